$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the data held in rows 4 and 5 (columns A, B, C, E — column D
# is "F" in both rows already, so it is left untouched).
$row4A = $ws.Range("A4").Value()
$row4B = $ws.Range("B4").Value()
$row4C = $ws.Range("C4").Value()
$row4E = $ws.Range("E4").Value()

$row5A = $ws.Range("A5").Value()
$row5B = $ws.Range("B5").Value()
$row5C = $ws.Range("C5").Value()
$row5E = $ws.Range("E5").Value()

$ws.Range("A4").Value = $row5A
$ws.Range("B4").Value = $row5B
$ws.Range("C4").Value = $row5C
$ws.Range("E4").Value = $row5E

$ws.Range("A5").Value = $row4A
$ws.Range("B5").Value = $row4B
$ws.Range("C5").Value = $row4C
$ws.Range("E5").Value = $row4E

# Move the active selection from B2 to C2.
$ws.Range("C2").Select()
